$d = $word.ActiveDocument

# 1. "follows his story" -> "follows Bob's story"
$d.Content.Find.Execute("follows his story", $true, $false, $false, $false, $false, $true, 1, $false, "follows Bob’s story", 2)

# 2. "don't pull out because" -> "don't withdrawal money because"
$d.Content.Find.Execute("don’t pull out because", $true, $false, $false, $false, $false, $true, 1, $false, "don’t withdrawal money because", 2)

# 3. "but rather to the progression" -> "but rather through the progression"
$d.Content.Find.Execute("but rather to the progression", $true, $false, $false, $false, $false, $true, 1, $false, "but rather through the progression", 2)

# 4. Insert ISSUE parenthetical after "...SP average values for the years mentioned. "
$r = $d.Content
$r.Find.Execute("SP average values for the years mentioned. ")
$r.Collapse(0)
$r.InsertAfter("(ISSUE: something in the last day or so and the button no longer clicks initially, please hover over first datapoint then click)")
$r2 = $d.Content
$r2.Find.Execute("ISSUE: ")
$r2.Bold = 1

# 5. "year end" -> "year-end"
$d.Content.Find.Execute("year end", $true, $false, $false, $false, $false, $true, 1, $false, "year-end", 2)

# 6/7. Fix double space and "scene by scene" -> "scene-by-scene"
$d.Content.Find.Execute("Altogether,  you", $true, $false, $false, $false, $false, $true, 1, $false, "Altogether, you", 2)
$d.Content.Find.Execute("scene by scene visual mapping", $true, $false, $false, $false, $false, $true, 1, $false, "scene-by-scene visual mapping", 2)

# 8. "visible and the axes" -> "visible, and the axes"
$d.Content.Find.Execute("stay visible and the axes", $true, $false, $false, $false, $false, $true, 1, $false, "stay visible, and the axes", 2)

# 9. "data,annotations" -> "data, annotations"
$d.Content.Find.Execute("data,annotations", $true, $false, $false, $false, $false, $true, 1, $false, "data, annotations", 2)

# 10. "not overwhelm the user." -> "not overwhelm the viewer."
$d.Content.Find.Execute("not overwhelm the user.", $true, $false, $false, $false, $false, $true, 1, $false, "not overwhelm the viewer.", 2)

# 11. "cleanest of my options as the wavy lines felt" -> "cleanest options as the wavy line boxes felt"
$d.Content.Find.Execute("cleanest of my options as the wavy lines felt", $true, $false, $false, $false, $false, $true, 1, $false, "cleanest options as the wavy line boxes felt", 2)

# 12. quotes around buttonclicks + "data/info is shown" -> "data/info and states are shown"
$d.Content.Find.Execute("are buttonclicks and slider indexes and define what data/info is shown on the page.", $true, $false, $false, $false, $false, $true, 1, $false, "are “buttonclicks” and slider indexes and define what data/info and states are shown on the page.", 2)
